# Auto-generated edit script applying scheduled price-refresh updates
# to the Goblin_Profits workbook (columns H/I/J/K/L/M/N per Leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 50000908
$ws.Range("I6").Value = 50000908
$ws.Range("K6").Value = 150002724
$ws.Range("M6").Value = -150002612

# Row 28
$ws.Range("H28").Value = 4606.885
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Row 64
$ws.Range("H64").Value = 8833.291999999999
$ws.Range("I64").Value = 4875
$ws.Range("K64").Value = 4875
$ws.Range("M64").Value = -4627

# Row 67
$ws.Range("H67").Value = 8833.291999999999
$ws.Range("I67").Value = 4875
$ws.Range("K67").Value = 4875
$ws.Range("M67").Value = -4017

# Row 88
$ws.Range("H88").Value = 9510.556
$ws.Range("J88").Value = 9510.556
$ws.Range("L88").Value = 9510.556
$ws.Range("N88").Value = -10322.556

# Row 91
$ws.Range("H91").Value = 9510.556
$ws.Range("J91").Value = 9510.556
$ws.Range("L91").Value = 9510.556
$ws.Range("N91").Value = -12318.556

# Row 113
$ws.Range("H113").Value = 3954.3
$ws.Range("J113").Value = 4799.8
$ws.Range("L113").Value = 4799.8
$ws.Range("N113").Value = -11307.8

# Row 116
$ws.Range("H116").Value = 3897
$ws.Range("I116").Value = 4726.3335
$ws.Range("K116").Value = 4726.3335
$ws.Range("M116").Value = -1284.3335

# Row 131
$ws.Range("H131").Value = 4063.2273
$ws.Range("I131").Value = 2024.5
$ws.Range("J131").Value = 9499.833000000001
$ws.Range("K131").Value = 6073.5
$ws.Range("L131").Value = 28499.499
$ws.Range("M131").Value = -1033.5
$ws.Range("N131").Value = -38579.499

# Row 132
$ws.Range("H132").Value = 2217.125
$ws.Range("I132").Value = 2287.4348
$ws.Range("J132").Value = 600
$ws.Range("K132").Value = 6862.3044
$ws.Range("L132").Value = 1800
$ws.Range("M132").Value = -4332.3044
$ws.Range("N132").Value = -6860

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 137
$ws.Range("H137").Value = 9867.65
$ws.Range("I137").Value = 14238.583
$ws.Range("J137").Value = 3311.25
$ws.Range("K137").Value = 42715.749
$ws.Range("L137").Value = 9933.75
$ws.Range("M137").Value = -40165.749
$ws.Range("N137").Value = -15033.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2709.5964
$ws.Range("I32").Value = 2624.9434
$ws.Range("K32").Value = 2624.9434
$ws.Range("M32").Value = -2337.9434

# Row 74
$ws.Range("H74").Value = 3424.76
$ws.Range("I74").Value = 3623.611
$ws.Range("K74").Value = 3623.611
$ws.Range("M74").Value = -2749.611

# Row 77
$ws.Range("H77").Value = 3424.76
$ws.Range("I77").Value = 3623.611
$ws.Range("K77").Value = 18118.055
$ws.Range("M77").Value = -13750.055

# Row 102
$ws.Range("H102").Value = 5492.533
$ws.Range("I102").Value = 2673.5
$ws.Range("K102").Value = 2673.5
$ws.Range("M102").Value = -1051.5

# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 122
$ws.Range("H122").Value = 5850866.5
$ws.Range("I122").Value = 6538720.5
$ws.Range("K122").Value = 19616161.5
$ws.Range("M122").Value = -19613711.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1683.3334
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 1820
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 1820
$ws.Range("M20").Value = -753
$ws.Range("N20").Value = -2314

# Row 94
$ws.Range("H94").Value = 112731.35
$ws.Range("I94").Value = 145365.8
$ws.Range("K94").Value = 145365.8
$ws.Range("M94").Value = -144914.8

# Row 105
$ws.Range("H105").Value = 2033.2778
$ws.Range("J105").Value = 2133
$ws.Range("L105").Value = 2133
$ws.Range("N105").Value = -5627

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 246
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 257.5
$ws.Range("K2").Value = 200
$ws.Range("L2").Value = 257.5
$ws.Range("M2").Value = -87
$ws.Range("N2").Value = -483.5

# Row 75
$ws.Range("H75").Value = 70000
$ws.Range("J75").Value = 70000
$ws.Range("L75").Value = 70000
$ws.Range("N75").Value = -71996

# Row 78
$ws.Range("H78").Value = 70000
$ws.Range("J78").Value = 70000
$ws.Range("L78").Value = 210000
$ws.Range("N78").Value = -219984

# Row 105
$ws.Range("H105").Value = 4750
$ws.Range("J105").Value = 4750
$ws.Range("L105").Value = 4750
$ws.Range("N105").Value = -8244

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2082.3845
$ws.Range("I68").Value = 2333
$ws.Range("J68").Value = 2007.2
$ws.Range("K68").Value = 6999
$ws.Range("L68").Value = 6021.6
$ws.Range("M68").Value = -6188
$ws.Range("N68").Value = -7643.6

# Row 71
$ws.Range("H71").Value = 2082.3845
$ws.Range("I71").Value = 2333
$ws.Range("J71").Value = 2007.2
$ws.Range("K71").Value = 20997
$ws.Range("L71").Value = 18064.8
$ws.Range("M71").Value = -16941
$ws.Range("N71").Value = -26176.8

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 4737.1787
$ws.Range("I97").Value = 1318.7391
$ws.Range("K97").Value = 1318.7391
$ws.Range("M97").Value = -822.7391

# Row 102
$ws.Range("H102").Value = 813.2857
$ws.Range("I102").Value = 722
$ws.Range("K102").Value = 722
$ws.Range("M102").Value = 900

# Row 113
$ws.Range("H113").Value = 38468320
$ws.Range("I113").Value = 76926024
$ws.Range("J113").Value = 10615.308
$ws.Range("K113").Value = 76926024
$ws.Range("L113").Value = 10615.308
$ws.Range("M113").Value = -76923854
$ws.Range("N113").Value = -14955.308

# Row 122
$ws.Range("H122").Value = 5694
$ws.Range("I122").Value = 5694
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17082
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14632
$ws.Range("N122").ClearContents()

# Row 123
$ws.Range("H123").Value = 69999.25
$ws.Range("J123").Value = 69999.25
$ws.Range("L123").Value = 69999.25
$ws.Range("N123").Value = -74899.25

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2867.7778
$ws.Range("I93").Value = 1519.8572
$ws.Range("J93").Value = 4319.385
$ws.Range("K93").Value = 1519.8572
$ws.Range("L93").Value = 4319.385
$ws.Range("M93").Value = -271.8571999999999
$ws.Range("N93").Value = -6815.385

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10383.333
$ws.Range("J62").Value = 11800
$ws.Range("L62").Value = 11800
$ws.Range("N62").Value = -13048

# Row 65
$ws.Range("H65").Value = 10383.333
$ws.Range("J65").Value = 11800
$ws.Range("L65").Value = 59000
$ws.Range("N65").Value = -65240

# Row 96
$ws.Range("H96").Value = 6074.6875
$ws.Range("I96").Value = 6145
$ws.Range("J96").Value = 5920
$ws.Range("K96").Value = 6145
$ws.Range("L96").Value = 5920
$ws.Range("M96").Value = -4772
$ws.Range("N96").Value = -8666

# Row 109
$ws.Range("H109").Value = 222000
$ws.Range("J109").Value = 222000
$ws.Range("L109").Value = 222000
$ws.Range("N109").Value = -224774
